$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 8033.357
$ws.Range("I15").Value = 8033.357
$ws.Range("K15").Value = 24100.071
$ws.Range("M15").Value = -23931.071
# row 29
$ws.Range("H29").Value = 1297.3334
$ws.Range("J29").Value = 2000
$ws.Range("L29").Value = 6000
$ws.Range("N29").Value = -6562
# row 116
$ws.Range("H116").Value = 20982.666
$ws.Range("J116").Value = 8946.25
$ws.Range("L116").Value = 8946.25
$ws.Range("N116").Value = -15830.25
# row 137
$ws.Range("H137").Value = 21725.533
$ws.Range("I137").Value = 35143.777
$ws.Range("J137").Value = 1598.1666
$ws.Range("K137").Value = 105431.331
$ws.Range("L137").Value = 4794.4998
$ws.Range("M137").Value = -102881.331
$ws.Range("N137").Value = -9894.4998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 27267.385
$ws.Range("I32").Value = 27906
$ws.Range("K32").Value = 27906
$ws.Range("M32").Value = -27619
# row 61
$ws.Range("H61").Value = 9623.615
$ws.Range("I61").Value = 999
$ws.Range("J61").Value = 19685.666
$ws.Range("K61").Value = 999
$ws.Range("L61").Value = 19685.666
$ws.Range("M61").Value = -787
$ws.Range("N61").Value = -20109.666
# row 74
$ws.Range("H74").Value = 382902
$ws.Range("I74").Value = 750652.6
$ws.Range("K74").Value = 750652.6
$ws.Range("M74").Value = -749778.6
# row 77
$ws.Range("H77").Value = 382902
$ws.Range("I77").Value = 750652.6
$ws.Range("K77").Value = 3753263
$ws.Range("M77").Value = -3748895
# row 122
$ws.Range("H122").Value = 3240.625
$ws.Range("I122").Value = 2782.4
$ws.Range("K122").Value = 8347.200000000001
$ws.Range("M122").Value = -5897.200000000001
# row 132
$ws.Range("H132").Value = 2530.3076
$ws.Range("I132").Value = 2089.4
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 6268.200000000001
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3738.200000000001
$ws.Range("N132").Value = -17060
# row 136
$ws.Range("H136").Value = 9623.615
$ws.Range("I136").Value = 999
$ws.Range("J136").Value = 19685.666
$ws.Range("K136").Value = 2997
$ws.Range("L136").Value = 59056.99800000001
$ws.Range("M136").Value = -447
$ws.Range("N136").Value = -64156.99800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 3637.111
$ws.Range("I134").Value = 3637.111
$ws.Range("K134").Value = 10911.333
$ws.Range("M134").Value = -8376.332999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 3126349.2
$ws.Range("I31").Value = 3572221
$ws.Range("K31").Value = 3572221
$ws.Range("M31").Value = -3571926
# row 34
$ws.Range("H34").Value = 3126349.2
$ws.Range("I34").Value = 3572221
$ws.Range("K34").Value = 3572221
$ws.Range("M34").Value = -3572019
# row 94
$ws.Range("H94").Value = 3446.0833
$ws.Range("I94").Value = 2071.2
$ws.Range("J94").Value = 4428.143
$ws.Range("K94").Value = 2071.2
$ws.Range("L94").Value = 4428.143
$ws.Range("M94").Value = -1620.2
$ws.Range("N94").Value = -5330.143
# row 99
$ws.Range("H99").Value = 3660.923
$ws.Range("I99").Value = 2199.3333
$ws.Range("K99").Value = 2199.3333
$ws.Range("M99").Value = -701.3332999999998
# row 126
$ws.Range("H126").Value = 3660.923
$ws.Range("I126").Value = 2199.3333
$ws.Range("K126").Value = 6597.999899999999
$ws.Range("M126").Value = -4127.999899999999
# row 132
$ws.Range("H132").Value = 39883.96
$ws.Range("I132").Value = 53509.473
$ws.Range("J132").Value = 2900.4285
$ws.Range("K132").Value = 160528.419
$ws.Range("L132").Value = 8701.2855
$ws.Range("M132").Value = -157998.419
$ws.Range("N132").Value = -13761.2855
# row 134
$ws.Range("H134").Value = 3335.3
$ws.Range("I134").Value = 2865.0625
$ws.Range("J134").Value = 5216.25
$ws.Range("K134").Value = 8595.1875
$ws.Range("L134").Value = 15648.75
$ws.Range("M134").Value = -6060.1875
$ws.Range("N134").Value = -20718.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 11
$ws.Range("H11").Value = 100649
$ws.Range("I11").Value = 111721.11
$ws.Range("J11").Value = 1000
$ws.Range("K11").Value = 335163.33
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = -335023.33
$ws.Range("N11").Value = -3280
# row 131
$ws.Range("H131").Value = 166379.16
$ws.Range("J131").Value = 2041.5
$ws.Range("L131").Value = 6124.5
$ws.Range("N131").Value = -16204.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 59
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()
# row 70
$ws.Range("H70").Value = 7405.1055
$ws.Range("I70").Value = 8247.615
$ws.Range("J70").Value = 5579.6665
$ws.Range("K70").Value = 8247.615
$ws.Range("L70").Value = 5579.6665
$ws.Range("M70").Value = -7977.615
$ws.Range("N70").Value = -6119.6665
# row 73
$ws.Range("H73").Value = 7405.1055
$ws.Range("I73").Value = 8247.615
$ws.Range("J73").Value = 5579.6665
$ws.Range("K73").Value = 8247.615
$ws.Range("L73").Value = 5579.6665
$ws.Range("M73").Value = -7311.615
$ws.Range("N73").Value = -7451.6665
# row 113
$ws.Range("H113").Value = 1958.8
$ws.Range("I113").Value = 1949.75
$ws.Range("K113").Value = 1949.75
$ws.Range("M113").Value = 220.25
# row 122
$ws.Range("H122").Value = 5535.4287
$ws.Range("I122").Value = 5399.8
$ws.Range("J122").Value = 5874.5
$ws.Range("K122").Value = 16199.4
$ws.Range("L122").Value = 17623.5
$ws.Range("M122").Value = -13749.4
$ws.Range("N122").Value = -22523.5
# row 126
$ws.Range("H126").Value = 1795.1428
$ws.Range("J126").Value = 2183.5
$ws.Range("L126").Value = 6550.5
$ws.Range("N126").Value = -11490.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 68
$ws.Range("H68").Value = 3794.3333
$ws.Range("I68").Value = 3831.3333
$ws.Range("J68").Value = 3683.3333
$ws.Range("K68").Value = 3831.3333
$ws.Range("L68").Value = 3683.3333
$ws.Range("M68").Value = -3082.3333
$ws.Range("N68").Value = -5181.3333
# row 71
$ws.Range("H71").Value = 3794.3333
$ws.Range("I71").Value = 3831.3333
$ws.Range("J71").Value = 3683.3333
$ws.Range("K71").Value = 19156.6665
$ws.Range("L71").Value = 18416.6665
$ws.Range("M71").Value = -15412.6665
$ws.Range("N71").Value = -25904.6665
# row 132
$ws.Range("H132").Value = 3004.125
$ws.Range("I132").Value = 2645.147
$ws.Range("J132").Value = 3875.9285
$ws.Range("K132").Value = 7935.441
$ws.Range("L132").Value = 11627.7855
$ws.Range("M132").Value = -5405.441
$ws.Range("N132").Value = -16687.7855
# row 136
$ws.Range("H136").Value = 3683.0435
$ws.Range("I136").Value = 3164.8572
$ws.Range("K136").Value = 9494.571599999999
$ws.Range("M136").Value = -6944.571599999999
# row 140
$ws.Range("H140").Value = 130800
$ws.Range("J140").Value = 130800
$ws.Range("L140").Value = 130800
$ws.Range("N140").Value = -141160

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 100
$ws.Range("H100").Value = 645.6923
$ws.Range("I100").Value = 687.1111
$ws.Range("K100").Value = 1374.2222
$ws.Range("M100").Value = -833.2221999999999
# row 132
$ws.Range("H132").Value = 27378.928
$ws.Range("I132").Value = 36651.75
$ws.Range("K132").Value = 109955.25
$ws.Range("M132").Value = -107425.25
# row 136
$ws.Range("H136").Value = 18395.975
$ws.Range("I136").Value = 27557.916
$ws.Range("K136").Value = 82673.74800000001
$ws.Range("M136").Value = -80123.74800000001
